# Reorder the comma-separated "Arbitrage" tags in column G into a
# canonical priority order: CA, BF, CS (any other/unknown tags keep
# their relative order and are appended after the known ones).
#
# Only cells whose value actually contains a comma (i.e. multiple tags)
# are affected; single-tag cells are left untouched since their order
# is already trivially "correct".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$order = @("CA", "BF", "CS")

$used = $ws.UsedRange
$rowCount = $used.Rows.Count

for ($r = 1; $r -le $rowCount; $r++) {
    $cell = $ws.Cells.Item($r, 7)
    $val = $cell.Value2

    if ($val -ne $null -and $val -is [string] -and $val.Contains(",")) {
        $tokens = $val.Split(",")

        $newTokens = New-Object System.Collections.ArrayList

        # First append known tokens in canonical order.
        foreach ($known in $order) {
            foreach ($tok in $tokens) {
                if ($tok -eq $known) {
                    [void]$newTokens.Add($tok)
                }
            }
        }

        # Then append any tokens not part of the known set, preserving
        # their original relative order.
        foreach ($tok in $tokens) {
            if ($order -notcontains $tok) {
                [void]$newTokens.Add($tok)
            }
        }

        $newVal = [string]::Join(",", $newTokens)

        if ($newVal -ne $val) {
            $cell.Value2 = $newVal
        }
    }
}
